# edit.ps1 - Applies the "output folder feature" change:
#   - Input sheet: records the source .dat filename for the analysis (new D3 value + widened column D)
#   - Calc sheet: adds two new "Erfolgsrate" (success-rate) columns (BG/BH) with new header text,
#     re-derives several dependent values (229Th spike correction variant) and widens/narrows a few columns
#   - Results sheet: mirrors the updated dependent values and narrows column P
#   - Constants sheet: updates the R30_29 constant used in the recalculation

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Input"
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")

# New context cell: which raw Neptune data file this analysis line came from
$wsInput.Range("D3").Value = "Filename: C:\Neptune\User\Neptune\Data\UTh\2015\0815\011_7184.dat"

# Widen column D so the filename is readable (target raw width 66.7109375 chars;
# ColumnWidth only accepts 1/6-character increments, so use closest achievable value)
$wsInput.Columns.Item(4).ColumnWidth = 65.83333333333333

# ---------------------------------------------------------------------------
# Sheet "Calc"
# ---------------------------------------------------------------------------
$wsCalc = $wb.Worksheets.Item("Calc")

# New "Montefehler Erfolgsrate" (Monte-Carlo error success rate) columns
$wsCalc.Range("BG1").Value = "Unkorr. Montefehler Erfolgsrate"
$wsCalc.Range("BH1").Value = "Korr. Montefehler Erfolgsrate"
$wsCalc.Range("BG2").Value = "(%)"
$wsCalc.Range("BH2").Value = "(%)"
$wsCalc.Range("BG3").Value = 100
$wsCalc.Range("BH3").Value = 100

# Recalculated values (driven by updated R30_29 constant on the Constants sheet)
$wsCalc.Range("AP3").Value = 0.5304
$wsCalc.Range("AQ3").Value = 0.1878915514337181
$wsCalc.Range("AW3").Value = 0.53
$wsCalc.Range("AX3").Value = 0.5336101578840383
$wsCalc.Range("AY3").Value = 0.1877849769892541
$wsCalc.Range("BC3").Value = 0.5686648625490963
$wsCalc.Range("BE3").Value = 266.8050789420192
$wsCalc.Range("BF3").Value = 0.1890640966405403

# Column width tweaks (target raw widths shown in comments; ColumnWidth snaps to 1/6-char steps)
$wsCalc.Columns.Item(49).ColumnWidth = 8.833333333333334   # AW -> 9.7109375
$wsCalc.Columns.Item(51).ColumnWidth = 18.833333333333332  # AY -> 19.7109375
$wsCalc.Columns.Item(58).ColumnWidth = 19.833333333333332  # BF -> 20.7109375
$wsCalc.Columns.Item(59).ColumnWidth = 31.833333333333332  # BG (new) -> 32.7109375
$wsCalc.Columns.Item(60).ColumnWidth = 29.833333333333332  # BH (new) -> 30.7109375

# ---------------------------------------------------------------------------
# Sheet "Results"
# ---------------------------------------------------------------------------
$wsResults = $wb.Worksheets.Item("Results")

$wsResults.Range("N3").Value = 0.5304
$wsResults.Range("P3").Value = 0.53
$wsResults.Range("R3").Value = 0.5686648625490963

$wsResults.Columns.Item(16).ColumnWidth = 7.833333333333333  # P -> 8.7109375

# ---------------------------------------------------------------------------
# Sheet "Constants"
# ---------------------------------------------------------------------------
$wsConstants = $wb.Worksheets.Item("Constants")

# R30_29 constant updated from 4.8E-05 to 5E-05
$wsConstants.Range("B3").Value = 0.00005
